$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.981.44"
$ws.Range("E2").Value = "  +2.49%  "

$ws.Range("D3").Value = "3.603.40"
$ws.Range("E3").Value = "  +1.18%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "564.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.82%  "

$ws.Range("D7").Value = "3.601.22"
$ws.Range("E7").Value = "  +1.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.18%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.673"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.151"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000286"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("D15").Value = "4.187.10"
$ws.Range("E15").Value = "  +1.30%  "

$ws.Range("D16").Value = "3.610.63"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.85%  "

$ws.Range("D19").Value = "67.911.39"
$ws.Range("E19").Value = "  +2.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.86%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.15%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "676.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.425"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.65%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").Value = "0.0₃0767"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("D41").Value = "3.287.29"
$ws.Range("E41").Value = "  +9.19%  "

$ws.Range("E42").Value = "  +13.74%  "

$ws.Range("E43").Value = "  +4.26%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +31.53%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0418"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("E48").Value = "  +10.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("E50").Value = "  +0.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.15%  "

